$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update C31 value (0.7 -> 1, displayed as percentage)
$ws.Range("C31").Value = 1

# Row 37: add empty formatted cell C37 (same style as D11/C18/B33 placeholders)
$ws.Range("C37").Value = ""
$ws.Range("C37").Font.Underline = 2

# Row 38
$ws.Range("A38").Value = "hacer configurable la ip de la printer"
$ws.Range("B38").Value = "Lucas"
$ws.Range("C38").Value = "en proceso"

# Row 39
$ws.Range("A39").Value = "primer cuota - 1 mes mas"
$ws.Range("B39").Value = "Agustina"
$ws.Range("C39").Value = "en proceso"

# Row 40
$ws.Range("A40").Value = "sacar cartel de cliente asociado con éxito"
$ws.Range("B40").Value = "Agustina"
$ws.Range("C40").Value = 1
$ws.Range("C40").NumberFormat = "0%"

# Update selection to mimic final cursor position
$ws.Range("E40").Select()
